$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns to make room for the new metrics:
#  - "Average_GPCD_Change" goes right after "Average_PerCapita_._Change" (old column G)
#  - "ET_Changes" goes right after "ET_._Change" (old column H once shifted)
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()

# New column headers
$ws.Range("G1").Value = "Average_GPCD_Change"
$ws.Range("I1").Value = "ET_Changes"

# New "Average_GPCD_Change" values (column G, rows 2-7)
$ws.Range("G2").Value = 17
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 17
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 17
$ws.Range("G7").Value = 0

# New "ET_Changes" values (column I, rows 2-7)
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 9142
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 16237
$ws.Range("I6").Value = 0
$ws.Range("I7").Value = 22331

# Updated MI...AG values in column A for the HighET scenarios
$ws.Range("A3").Value = 557841.74555641
$ws.Range("A5").Value = 711807.63214346
$ws.Range("A7").Value = 846445.840554488
